$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F4").Value = 1
$ws.Range("F5").Value = 1
$ws.Range("F12").Value = 1
$ws.Range("F17").Value = 1
$ws.Range("F21").Value = $null
$ws.Range("F23").Value = $null
$ws.Range("G23").Formula = "=IF(F23=1,D23,0)"
